$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were dropped from the dataset: "RM 232" and "SC 92".
# Deleting row 26 ("RM 232") shifts everything up, so the old row 28 ("SC 92")
# becomes row 27; deleting row 27 next removes it too.
$ws.Range("A26").EntireRow.Delete()
$ws.Range("A27").EntireRow.Delete()

# Update column E (missing-data pattern changed: some filled in, some newly blanked)
$ws.Range("E2").Value = -7.2
$ws.Range("E3").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("E11").Value = -7.9
$ws.Range("E13").Value = ""
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("E25").Value = ""

# Row 29 is now "SC 119" (B becomes missing)
$ws.Range("B29").Value = ""

# Row 33 is now "SC 232" (B and E become filled in)
$ws.Range("B33").Value = -19.5
$ws.Range("E33").Value = -10.7
